$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scenario block cell references from row 33 to row 36
$ws.Range("D5").Value = "A36"
$ws.Range("D6").Value = "B36"
$ws.Range("D7").Value = "C36"
$ws.Range("D8").Value = "G36"
$ws.Range("D9").Value = "H36"
$ws.Range("D10").Value = "I36"
$ws.Range("D11").Value = "J36"

# Update selection/view to the newly edited block
$ws.Range("D5:D11").Select()
